$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46: FISV/Fiserv -> FTNT/Fortinet (Sector unchanged, Sub-Sector -> Systems Software)
$ws.Range("A46").Value = "FTNT"
$ws.Range("B46").Value = "Fortinet"
$ws.Range("D46").Value = "Systems Software"

# Row 47: FTNT/Fortinet -> GEHC/GE HealthCare (Sector -> Health Care, Sub-Sector -> Health Care Technology)
$ws.Range("A47").Value = "GEHC"
$ws.Range("B47").Value = "GE HealthCare"
$ws.Range("C47").Value = "Health Care"
$ws.Range("D47").Value = "Health Care Technology"
